$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1390.2
$ws.Range("J19").Value = 1218.5714
$ws.Range("L19").Value = 1218.5714
$ws.Range("N19").Value = -1568.5714
$ws.Range("H40").Value = 3214.1428
$ws.Range("I40").Value = 2415.4285
$ws.Range("K40").Value = 2415.4285
$ws.Range("M40").Value = -2240.4285
$ws.Range("H86").Value = 2395742
$ws.Range("I86").Value = 3155.0833
$ws.Range("J86").Value = 5266846.5
$ws.Range("K86").Value = 3155.0833
$ws.Range("L86").Value = 5266846.5
$ws.Range("M86").Value = -2032.0833
$ws.Range("N86").Value = -5269092.5
$ws.Range("H89").Value = 2395742
$ws.Range("I89").Value = 3155.0833
$ws.Range("J89").Value = 5266846.5
$ws.Range("K89").Value = 15775.4165
$ws.Range("L89").Value = 26334232.5
$ws.Range("M89").Value = -10159.4165
$ws.Range("N89").Value = -26345464.5
$ws.Range("H100").Value = 3994.5
$ws.Range("I100").Value = 2792.389
$ws.Range("J100").Value = 5797.6665
$ws.Range("K100").Value = 2792.389
$ws.Range("L100").Value = 5797.6665
$ws.Range("M100").Value = -2251.389
$ws.Range("N100").Value = -6879.6665
$ws.Range("H132").Value = 17090.773
$ws.Range("J132").Value = 500000
$ws.Range("L132").Value = 1500000
$ws.Range("N132").Value = -1505060
$ws.Range("H138").Value = 2281.3076
$ws.Range("J138").Value = 3590.6365
$ws.Range("L138").Value = 10771.9095
$ws.Range("N138").Value = -21051.9095

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4800.625
$ws.Range("I2").Value = 6001.4
$ws.Range("J2").Value = 2799.3333
$ws.Range("K2").Value = 6001.4
$ws.Range("L2").Value = 2799.3333
$ws.Range("M2").Value = -5888.4
$ws.Range("N2").Value = -3025.3333
$ws.Range("H4").Value = 100.333336
$ws.Range("I4").Value = 124
$ws.Range("J4").Value = 53
$ws.Range("K4").Value = 124
$ws.Range("L4").Value = 53
$ws.Range("M4").Value = -8
$ws.Range("N4").Value = -285
$ws.Range("H5").Value = 62.5
$ws.Range("I5").Value = 73.333336
$ws.Range("J5").Value = 30
$ws.Range("K5").Value = 73.333336
$ws.Range("L5").Value = 30
$ws.Range("M5").Value = 38.666664
$ws.Range("N5").Value = -254
$ws.Range("H32").Value = 55568300
$ws.Range("I32").Value = 55568300
$ws.Range("K32").Value = 55568300
$ws.Range("M32").Value = -55568013
$ws.Range("H45").Value = 1637.5
$ws.Range("I45").Value = 1637.5
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1637.5
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1260.5
$ws.Range("N45").ClearContents()
$ws.Range("H74").Value = 1769.76
$ws.Range("I74").Value = 1756.7727
$ws.Range("K74").Value = 1756.7727
$ws.Range("M74").Value = -882.7727
$ws.Range("H77").Value = 1769.76
$ws.Range("I77").Value = 1756.7727
$ws.Range("K77").Value = 8783.863499999999
$ws.Range("M77").Value = -4415.863499999999
$ws.Range("H116").Value = 4800.625
$ws.Range("I116").Value = 6001.4
$ws.Range("J116").Value = 2799.3333
$ws.Range("K116").Value = 6001.4
$ws.Range("L116").Value = 2799.3333
$ws.Range("M116").Value = -3707.4
$ws.Range("N116").Value = -7387.3333
$ws.Range("H122").Value = 4336.1055
$ws.Range("I122").Value = 2540.75
$ws.Range("J122").Value = 7413.857
$ws.Range("K122").Value = 7622.25
$ws.Range("L122").Value = 22241.571
$ws.Range("M122").Value = -5172.25
$ws.Range("N122").Value = -27141.571

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4800.625
$ws.Range("I3").Value = 6001.4
$ws.Range("J3").Value = 2799.3333
$ws.Range("K3").Value = 6001.4
$ws.Range("L3").Value = 2799.3333
$ws.Range("M3").Value = -5887.4
$ws.Range("N3").Value = -3027.3333
$ws.Range("H4").Value = 62.5
$ws.Range("I4").Value = 73.333336
$ws.Range("J4").Value = 30
$ws.Range("K4").Value = 73.333336
$ws.Range("L4").Value = 30
$ws.Range("M4").Value = 41.666664
$ws.Range("N4").Value = -260
$ws.Range("H20").Value = 1142
$ws.Range("I20").Value = 968.4286
$ws.Range("J20").Value = 1749.5
$ws.Range("K20").Value = 968.4286
$ws.Range("L20").Value = 1749.5
$ws.Range("M20").Value = -721.4286
$ws.Range("N20").Value = -2243.5
$ws.Range("H22").Value = 3781.125
$ws.Range("I22").Value = 4458.1665
$ws.Range("J22").Value = 1750
$ws.Range("K22").Value = 4458.1665
$ws.Range("L22").Value = 1750
$ws.Range("M22").Value = -4285.1665
$ws.Range("N22").Value = -2096
$ws.Range("H134").Value = 6287.048
$ws.Range("I134").Value = 5686.364
$ws.Range("J134").Value = 6947.8
$ws.Range("K134").Value = 17059.092
$ws.Range("L134").Value = 20843.4
$ws.Range("M134").Value = -14524.092
$ws.Range("N134").Value = -25913.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 137.45833
$ws.Range("I7").Value = 26.7
$ws.Range("K7").Value = 26.7
$ws.Range("M7").Value = 86.3
$ws.Range("H12").Value = 700
$ws.Range("I12").Value = 400
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 400
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = -230
$ws.Range("N12").Value = -1340
$ws.Range("H132").Value = 3658.4285
$ws.Range("I132").Value = 1927.5
$ws.Range("K132").Value = 5782.5
$ws.Range("M132").Value = -3252.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 889.95
$ws.Range("I113").Value = 665.3333
$ws.Range("K113").Value = 1995.9999
$ws.Range("M113").Value = 174.0001
$ws.Range("H121").Value = 1330.6364
$ws.Range("I121").Value = 457.72726
$ws.Range("K121").Value = 1373.18178
$ws.Range("M121").Value = -63.18177999999989
$ws.Range("H122").Value = 962.2857
$ws.Range("J122").Value = 1040.8889
$ws.Range("L122").Value = 9368.000099999999
$ws.Range("N122").Value = -14268.0001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 15875358
$ws.Range("J80").Value = 25643010
$ws.Range("L80").Value = 25643010
$ws.Range("N80").Value = -25645006
$ws.Range("H83").Value = 15875358
$ws.Range("J83").Value = 25643010
$ws.Range("L83").Value = 128215050
$ws.Range("N83").Value = -128225034
$ws.Range("H122").Value = 3125.4285
$ws.Range("I122").Value = 2038.2142
$ws.Range("J122").Value = 5299.857
$ws.Range("K122").Value = 6114.642599999999
$ws.Range("L122").Value = 15899.571
$ws.Range("M122").Value = -3664.642599999999
$ws.Range("N122").Value = -20799.571
$ws.Range("H126").Value = 5933
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 5933
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 17799
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -22739

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 135.92857
$ws.Range("I16").Value = 138.6923
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 138.6923
$ws.Range("L16").Value = 100
$ws.Range("M16").Value = 31.30770000000001
$ws.Range("N16").Value = -440
$ws.Range("H116").Value = 189999.5
$ws.Range("J116").Value = 189999.5
$ws.Range("L116").Value = 189999.5
$ws.Range("N116").Value = -199177.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 348868.75
$ws.Range("I122").Value = 529361.3
$ws.Range("K122").Value = 1588083.9
$ws.Range("M122").Value = -1585633.9
$ws.Range("H126").Value = 2094.889
$ws.Range("I126").Value = 2106.125
$ws.Range("K126").Value = 6318.375
$ws.Range("M126").Value = -3848.375
